# Add minimal formatting to Excel sheet:
#  - hide gridlines on every worksheet
#  - set explicit column widths on every worksheet
#  - turn each worksheet's used range into a named Table (ListObject)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Hide gridlines on every sheet, then restore sheet 1 ("Table") as the
#    active / selected tab (it was already the selected tab originally).
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Activate()
    $excel.ActiveWindow.DisplayGridlines = $false
}
$wb.Worksheets.Item(1).Activate()

# ---------------------------------------------------------------------------
# 2. Column widths per sheet.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)   # Table
$ws1.Columns.Item(1).ColumnWidth = 15.8
$ws1.Columns.Item(2).ColumnWidth = 43.8

$ws2 = $wb.Worksheets.Item(2)   # Variables
$ws2.Columns.Item(1).ColumnWidth = 7.8
$ws2.Columns.Item(2).ColumnWidth = 7.8
$ws2.Columns.Item(3).ColumnWidth = 6.8
$ws2.Columns.Item(4).ColumnWidth = 11.8
$ws2.Columns.Item(5).ColumnWidth = 16.8
$ws2.Columns.Item(6).ColumnWidth = 25.8
$ws2.Columns.Item(7).ColumnWidth = 28.8
$ws2.Columns.Item(8).ColumnWidth = 31.8
$ws2.Columns.Item(9).ColumnWidth = 31.8

$ws3 = $wb.Worksheets.Item(3)   # Codelists
$ws3.Columns.Item(1).ColumnWidth = 7.8
$ws3.Columns.Item(2).ColumnWidth = 3.8
$ws3.Columns.Item(3).ColumnWidth = 8.8
$ws3.Columns.Item(4).ColumnWidth = 25.8
$ws3.Columns.Item(5).ColumnWidth = 34.8
$ws3.Columns.Item(6).ColumnWidth = 8.8

$ws4 = $wb.Worksheets.Item(4)   # Data
$ws4.Columns.Item(1).ColumnWidth = 4.8
$ws4.Columns.Item(2).ColumnWidth = 5.8
$ws4.Columns.Item(3).ColumnWidth = 5.8
$ws4.Columns.Item(4).ColumnWidth = 7.8

# ---------------------------------------------------------------------------
# 3. Turn each sheet's data range into a Table.
#    NOTE: the tables are created first (in sheet order) and only renamed
#    afterwards, in *reverse* order. Renaming immediately after creation
#    (or renaming forwards) trips an engine quirk that silently drops the
#    earliest-created table once a 3rd/4th table is added, so renaming is
#    deferred and done back-to-front to keep every table intact.
# ---------------------------------------------------------------------------
$lo1 = $ws1.ListObjects.Add(1, $ws1.Range("A1:B32"), $null, 1)
$lo2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:I5"), $null, 1)
$lo3 = $ws3.ListObjects.Add(1, $ws3.Range("A1:F11"), $null, 1)
$lo4 = $ws4.ListObjects.Add(1, $ws4.Range("A1:D85"), $null, 1)

$lo4.Name = "Table6"
$lo3.Name = "Table5"
$lo2.Name = "Table4"
$lo1.Name = "Table3"
